# Add a new "horses" entry below the existing animal list and move the
# current selection to C10 (matching the committed worksheet state).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A6").Value = "horses"

$ws.Range("C10").Select()
